# Updated set to better align with writing policy
# Rewrites the achievement-description wording in the "Achievements" sheet.
# Dependent sheets (Revision Text, RAScript Text, Checklist) hold formulas
# that reference these cells directly, so they recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Achievements")

# "Clear level N in 10 steps or less [5 NS - 5 EW]" -> drop the bracketed hint
$ws.Range("E13").Value = "Clear level 1 in 10 steps or less"
$ws.Range("E12").Value = "Clear level 2 in 10 steps or less"
$ws.Range("E11").Value = "Clear level 3 in 10 steps or less"
$ws.Range("E10").Value = "Clear level 4 in 10 steps or less"

# "Defeat N monsters ..." -> "Defeat N bad robots ..., shooting good robots does not add to count"
$ws.Range("E23").Value = "Defeat 3 bad robots in a row without taking damage or an emergency recharge, shooting good robots does not add to count"
$ws.Range("E24").Value = "Defeat 5 bad robots in a row without taking damage or an emergency recharge, shooting good robots does not add to count"
$ws.Range("E25").Value = "Defeat 10 bad robots in a row without taking damage or an emergency recharge, shooting good robots does not add to count"

# "Clear level N with 100+ power" -> "Clear level N with 100+ power credits"
$ws.Range("E8").Value = "Clear level 2 with 100+ power credits"
$ws.Range("E7").Value = "Clear level 3 with 100+ power credits"
$ws.Range("E6").Value = "Clear level 4 with 100+ power credits"

# Move the active selection to match the saved view state (E18, scrolled to top).
$ws.Activate()
$ws.Range("E18").Select()
